$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings are preserved exactly as text
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.003.33"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "2.337.83"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "306.38"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").Value = "101.28"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  -4.35%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("D10").Value = "34.92"
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("D11").Value = "52.26"
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").Value = "0.0801"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").Value = "15.87"
$ws.Range("E15").Value = "  +5.70%  "
$ws.Range("D16").Value = "2.343.83"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "0.816"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "42.936.58"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0913"
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "11.76"
$ws.Range("E21").Value = "  -4.64%  "
$ws.Range("D22").Value = "67.93"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "237.13"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "25.46"
$ws.Range("E27").Value = "  +3.42%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.32"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "35.04"
$ws.Range("E29").Value = "  -5.11%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "9.43"
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "161.09"
$ws.Range("E31").Value = "  -4.58%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.13"
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "17.61"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").Value = "  +5.59%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.0729"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.86"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").Value = "2.92"
$ws.Range("E39").Value = "  -4.90%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.113"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "2.47"
$ws.Range("E42").Value = "  +6.72%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.008.81"
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0287"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "18.82"
$ws.Range("E45").Value = "  -4.29%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "10.20"
$ws.Range("E46").Value = "  +3.36%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").Value = "55.85"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "2.89"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.562.53"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "4.73"
$ws.Range("E51").Value = "  +3.20%  "
